# Apply the changes described by the commit:
#  - Framework for test suite started: add a new task row (T5.0) to the
#    Sprint1 sheet, and set Status (DC/NO) values for a couple Backlog rows.
#  - Minor updates to team report (active sheet / selection changes).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Backlog sheet: add "Status" values in column D for rows 2-5
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("D2").Value = "DC"
$backlog.Range("D3").Value = "DC"
$backlog.Range("D4").Value = "NO"
$backlog.Range("D5").Value = "NO"

# Update the sheet's selection / scroll position
$backlog.Range("D5").Select()

# ---------------------------------------------------------------------
# Sprint1 sheet: add a new task row (T5.0) for the testing suite
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Range("A14").Value = "T5.0"
$sprint1.Range("B14").Value = "begin developing automated testing suite"
$sprint1.Range("C14").Value = "DC/NO"
$sprint1.Range("D14").Value = "in progress"
$sprint1.Range("E14").Value = 75
$sprint1.Range("F14").Value = 75

# Update the sheet's selection
$sprint1.Range("A14").Select()

# ---------------------------------------------------------------------
# Workbook view: switch the active/selected tab from Sprint1 (index 5)
# to Backlog (index 2), and make Backlog the tab-selected sheet.
# ---------------------------------------------------------------------
$backlog.Activate()
